# Auto-generated edit script: updates the cryptos price/volume table
# (mirrors a scheduled "Updated cryptos list ... with GitHub Actions" data refresh).
#
# All Price (D) / Volume(1h) (E) cells in this sheet are stored as literal text
# (not numbers), and a couple of rows also swap coin identity (name/link/price/volume).
# For D-column values that look like plain decimals, Excel would normally auto-detect
# them as numbers on assignment, so we briefly force Text format, write the value, then
# clear the formatting again so the cell keeps its original (default) style but retains
# the text value exactly as in the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.133.12'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.426.83'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.12'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.23'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.19%  '
$ws.Range('D9').Value = '2.428.04'
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.75%  '
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.32'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.70%  '
$ws.Range('E15').Value = '  -1.56%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = '61.996.86'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').Value = '2.430.42'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.28'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.83'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.85'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.15'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.46'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.93%  '
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.61'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '559.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('D28').Value = '2.545.84'
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '0.0₃0938'
$ws.Range('E30').Value = '  -0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.25'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('E33').Value = '  -2.46%  '
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('E35').Value = '  -2.07%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.78'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  -0.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.52'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '151.25'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.71'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.81'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.27'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.85'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0531'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.01'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.596'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E51').Value = '  +0.06%  '
